# Updates the Poland iii-liga-group-iv 2023-2024 odds sheet:
#  1) Re-orders several same-kickoff-date fixtures (rows whose F:V content
#     needs to be rotated among rows that share the same match date) so the
#     row order matches the freshly re-scraped source ordering.
#  2) Appends 4 newly scraped fixtures (25/26-11-2023) as rows 138-141.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: re-order the same-date fixture rows -------------------------
# Map of targetRow -> sourceRow: the new content of targetRow (columns F:V)
# is the OLD content that currently sits in sourceRow. All sources are read
# into memory before any writes happen so the (disjoint) row cycles don't
# clobber each other mid-flight.
$rowMap = @{
    28 = 30; 30 = 31; 31 = 28;
    65 = 68; 67 = 65; 68 = 67;
    73 = 75; 74 = 73; 75 = 74;
    96 = 97; 97 = 96;
    103 = 104; 104 = 103;
    105 = 106; 106 = 107; 107 = 105;
    112 = 114; 114 = 112;
    121 = 122; 122 = 121;
    132 = 133; 133 = 132
}

$snapshot = @{}
foreach ($key in $rowMap.Keys) {
    $srcRow = $rowMap[$key]
    if (-not $snapshot.ContainsKey($srcRow)) {
        $snapshot[$srcRow] = $ws.Range("F" + $srcRow + ":V" + $srcRow).Value()
    }
}

foreach ($key in $rowMap.Keys) {
    $srcRow = $rowMap[$key]
    $ws.Range("F" + $key + ":V" + $key).Value = $snapshot[$srcRow]
}

# --- Step 2: append the 4 newly scraped fixtures --------------------------
$newRows = @(
    @{ Row=138; A=137; E=45255.5;             F="KS Wieczysta Krakow"; G=5; H="Karpaty Krosno";       I=0;
       J=1.05; K="25/11/2023 00:12"; L=1.07;  M="25/11/2023 11:09";
       N=12.24; O="25/11/2023 00:12"; P=11.96; Q="25/11/2023 11:46";
       R=14.55; S="25/11/2023 00:12"; T=13.99; U="25/11/2023 11:46";
       V="https://www.betexplorer.com/football/poland/iii-liga-group-iv/ks-wieczysta-krakow-ks-karpaty-krosno/b9XbU8sm/" },

    @{ Row=139; A=138; E=45255.52083333334;   F="Unia Tarnow";         G=2; H="Garbarnia";             I=1;
       J=2.4;  K="25/11/2023 01:42"; L=2.53;  M="25/11/2023 12:23";
       N=3.46; O="25/11/2023 01:42"; P=3.43;  Q="25/11/2023 12:24";
       R=2.4;  S="25/11/2023 01:42"; T=2.39;  U="25/11/2023 12:24";
       V="https://www.betexplorer.com/football/poland/iii-liga-group-iv/unia-tarnow-garbarnia/46z7SnCa/" },

    @{ Row=140; A=139; E=45255.54166666666;   F="Wisloka Debica";      G=1; H="Wislanie Jaskowice";    I=3;
       J=2.24; K="25/11/2023 02:12"; L=2.43;  M="25/11/2023 11:44";
       N=3.39; O="25/11/2023 02:12"; P=3.29;  Q="25/11/2023 11:44";
       R=2.63; S="25/11/2023 02:12"; T=2.57;  U="25/11/2023 11:44";
       V="https://www.betexplorer.com/football/poland/iii-liga-group-iv/wisloka-debica-wislanie-jaskowice/hWW2TSdg/" },

    @{ Row=141; A=140; E=45256.5;             F="Sokol Sieniawa";      G=1; H="Siarka Tarnobrzeg";     I=5;
       J=4.4;  K="26/11/2023 01:12"; L=5.5;   M="26/11/2023 11:55";
       N=4.27; O="26/11/2023 01:12"; P=4.56;  Q="26/11/2023 11:55";
       R=1.51; S="26/11/2023 01:12"; T=1.42;  U="26/11/2023 11:55";
       V="https://www.betexplorer.com/football/poland/iii-liga-group-iv/sokol-sieniawa-siarka-tarnobrzeg/0YiOzB4J/" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value  = $nr.A          # Indice
    $ws.Cells.Item($r, 2).Value  = "poland"        # pais
    $ws.Cells.Item($r, 3).Value  = "iii-liga-group-iv"  # torneio
    $ws.Cells.Item($r, 4).Value  = "2023-2024"     # temporada
    $ws.Cells.Item($r, 5).Value  = $nr.E           # data_partida
    $ws.Cells.Item($r, 6).Value  = $nr.F           # home
    $ws.Cells.Item($r, 7).Value  = $nr.G           # home_ft_gols
    $ws.Cells.Item($r, 8).Value  = $nr.H           # away
    $ws.Cells.Item($r, 9).Value  = $nr.I           # away_ft_gols
    $ws.Cells.Item($r, 10).Value = $nr.J           # home_opening_odds
    $ws.Cells.Item($r, 11).Value = $nr.K           # home_opening_data_hora
    $ws.Cells.Item($r, 12).Value = $nr.L           # home_closing_odds
    $ws.Cells.Item($r, 13).Value = $nr.M           # home_closing_data_hora
    $ws.Cells.Item($r, 14).Value = $nr.N           # draw_opening_odds
    $ws.Cells.Item($r, 15).Value = $nr.O           # draw_opening_data_hora
    $ws.Cells.Item($r, 16).Value = $nr.P           # draw_closing_odds
    $ws.Cells.Item($r, 17).Value = $nr.Q           # draw_closing_data_hora
    $ws.Cells.Item($r, 18).Value = $nr.R           # away_opening_odds
    $ws.Cells.Item($r, 19).Value = $nr.S           # away_opening_data_hora
    $ws.Cells.Item($r, 20).Value = $nr.T           # away_closing_odds
    $ws.Cells.Item($r, 21).Value = $nr.U           # away_closing_data_hora
    $ws.Cells.Item($r, 22).Value = $nr.V           # url_partida
}

Write-Host "Done: reordered same-date fixtures and appended 4 new rows (138-141)."
